$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price (column D) cells to Text format before assignment so that
# numeric-looking strings (e.g. "82.80", "1.00") are not silently converted
# into numbers (which would drop significant trailing zeros).
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D8", "D9", "D12", "D13", "D14", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D32", "D34", "D35", "D37", "D38", "D39", "D41", "D43", "D44", "D47", "D49", "D50", "D51")
foreach ($ref in $priceCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.906.64"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "3.080.20"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "578.56"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "168.87"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.077.26"
$ws.Range("E8").Value = "  -1.25%  "
$ws.Range("D9").Value = "0.514"
$ws.Range("E9").Value = "  -1.49%  "
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("E11").Value = "  -2.39%  "
$ws.Range("D12").Value = "0.472"
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("D13").Value = "0.0000241"
$ws.Range("E13").Value = "  -2.70%  "
$ws.Range("D14").Value = "36.24"
$ws.Range("E14").Value = "  -2.81%  "
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("D16").Value = "3.589.40"
$ws.Range("D17").Value = "66.788.46"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").Value = "7.02"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "3.083.30"
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("D20").Value = "16.41"
$ws.Range("E20").Value = "  +0.85%  "
$ws.Range("D21").Value = "484.65"
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("D22").Value = "7.73"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "0.689"
$ws.Range("E23").Value = "  -3.38%  "
$ws.Range("D24").Value = "82.80"
$ws.Range("E24").Value = "  -1.31%  "
$ws.Range("D25").Value = "12.88"
$ws.Range("E25").Value = "  -3.30%  "
$ws.Range("E26").Value = "  -2.90%  "
$ws.Range("D27").Value = "10.25"
$ws.Range("E27").Value = "  +2.42%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -3.46%  "
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "27.85"
$ws.Range("E32").Value = "  -3.04%  "
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").Value = "0.0₃0912"
$ws.Range("E34").Value = "  -6.46%  "
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("D37").Value = "0.953"
$ws.Range("E37").Value = "  -2.56%  "
$ws.Range("D38").Value = "46.45"
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").Value = "0.123"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("E40").Value = "  -5.19%  "
$ws.Range("D41").Value = "0.302"
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("D43").Value = "2.771.75"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "373.88"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("D47").Value = "134.93"
$ws.Range("E47").Value = "  -0.73%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").Value = "24.41"
$ws.Range("E49").Value = "  -1.68%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "0.106"
$ws.Range("E50").Value = "  -1.36%  "
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "2.14"
$ws.Range("E51").Value = "  -3.04%  "

# Restore the default (Normal) style on the price cells so no stray number
# format / quote-prefix style index lingers on the saved cells.
foreach ($ref in $priceCells) {
    $ws.Range($ref).Style = "Normal"
}
